# Split the pycairo version "1.8.10" into three runs reading "1." + "10" + ".0"
# (i.e. the displayed version becomes "1.10.0"), matching the py2cairo-1.10.0
# download link referenced just below it.

$d = $word.ActiveDocument

# Locate the paragraph whose whole text is the version string "1.8.10".
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "1.8.10") {
        $target = $p.Range
        break
    }
}

$pStart = $target.Start

# Run 1: rewrite "1.8.10" down to "1." (keeps this as a single, first run).
$r1 = $d.Range($pStart, $pStart + 6)
$r1.Text = "1."

# Run 2: insert "10" right after it, in its own run (toggle Italic off/on so
# the engine doesn't silently re-merge it into the previous run).
$r2 = $d.Range($pStart + 2, $pStart + 2)
$r2.InsertAfter("10")
$r2run = $d.Range($pStart + 2, $pStart + 4)
$r2run.Font.Italic = $false
$r2run.Font.Italic = $true

# Run 3: insert ".0" right after that, again as its own run.
$r3 = $d.Range($pStart + 4, $pStart + 4)
$r3.InsertAfter(".0")
$r3run = $d.Range($pStart + 4, $pStart + 6)
$r3run.Font.Italic = $false
$r3run.Font.Italic = $true

Write-Output "pycairo version paragraph now reads: [$($target.Text)]"
